$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy the header style (bold, border, centered) from an existing header cell (F1) to G1:H1
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update existing MSE / R2 / MAE values and add Elapsed Time / CPU columns

# Row 2
$ws.Range("B2").Value = 0.5289725419958162
$ws.Range("C2").Value = 0.9894670333515272
$ws.Range("D2").Value = 0.6079239902858726
$ws.Range("G2").Value = 0.4788041146331428
$ws.Range("H2").Value = 0.997

# Row 3
$ws.Range("B3").Value = 0.2417608117777805
$ws.Range("C3").Value = 0.995274778853975
$ws.Range("D3").Value = 0.3930138500027158
$ws.Range("G3").Value = 0.4788041146331428
$ws.Range("H3").Value = 0.997

# Row 4
$ws.Range("B4").Value = 0.2772240630221857
$ws.Range("C4").Value = 0.994666825474385
$ws.Range("D4").Value = 0.4269450761411213
$ws.Range("G4").Value = 0.4788041146331428
$ws.Range("H4").Value = 0.997

# Row 5
$ws.Range("B5").Value = 0.4126945062468326
$ws.Range("C5").Value = 0.9918622293659692
$ws.Range("D5").Value = 0.494246395357917
$ws.Range("G5").Value = 0.4788041146331428
$ws.Range("H5").Value = 0.997
